# Replace the 2nd and 3rd "Ime Priimek clana N" placeholder paragraphs in the
# "Text Placeholder 8" shape on slide 1 with the real team member names
# "Luka Rus" and "Miha Markočič".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Find the members/text placeholder shape by name (falls back to the known
# index if the name ever changes).
$shp = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $cand = $s.Shapes.Item($i)
    if ($cand.Name -eq "Text Placeholder 8") {
        $shp = $cand
        break
    }
}
if ($shp -eq $null) {
    $shp = $s.Shapes.Item(3)
}

$tr = $shp.TextFrame.TextRange

# Paragraph 1 = "Žan Koren Kern" (untouched)
# Paragraph 2 = "Ime Priimek člana 2" -> "Luka Rus"
# Paragraph 3 = "Ime Priimek člana 3" -> "Miha Markočič"
# Paragraph 4 = "Ime Priimek člana 4" (untouched)

$para2 = $tr.Paragraphs(2, 1)
$para2.Text = "Luka Rus"

$para3 = $tr.Paragraphs(3, 1)
$para3.Text = "Miha Markočič"
